$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.38225564643511
$ws.Range("C2").Value = 7.466367187507856
$ws.Range("D2").Value = 4.458777462262016
$ws.Range("E2").Value = 11.36521475653975
$ws.Range("F2").Value = 60.79824859264365
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 10.83153443283545
$ws.Range("K2").Value = 15.0308099310607

$ws.Range("B3").Value = 15.31761715088792
$ws.Range("C3").Value = 7.461854550307002
$ws.Range("D3").Value = 4.511138247076593
$ws.Range("E3").Value = 11.41148789700226
$ws.Range("F3").Value = 59.82695907145681
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 10.82488498585506
$ws.Range("K3").Value = 15.02419668091874

$ws.Range("B4").Value = 15.2838957846525
$ws.Range("C4").Value = 7.463283789944272
$ws.Range("D4").Value = 4.545934970893296
$ws.Range("E4").Value = 11.4435079774898
$ws.Range("F4").Value = 59.22646136410962
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 10.82240464489282
$ws.Range("K4").Value = 15.02566999335031

$ws.Range("B5").Value = 15.27166667622154
$ws.Range("C5").Value = 7.464919303483766
$ws.Range("D5").Value = 4.560773955365311
$ws.Range("E5").Value = 11.45746059880442
$ws.Range("F5").Value = 58.98092973084028
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 10.82179724639114
$ws.Range("K5").Value = 15.02766123609697

$ws.Range("B6").Value = 15.26972770708934
$ws.Range("C6").Value = 7.465254364239837
$ws.Range("D6").Value = 4.563277573708235
$ws.Range("E6").Value = 11.45983193716586
$ws.Range("F6").Value = 58.94011588263026
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 10.82172075899006
$ws.Range("K6").Value = 15.02807583895394

$ws.Range("B7").Value = 15.28372471993644
$ws.Range("C7").Value = 7.463301588364041
$ws.Range("D7").Value = 4.54613243477763
$ws.Range("E7").Value = 11.44369249070263
$ws.Range("F7").Value = 59.2231530924363
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 10.82239481969393
$ws.Range("K7").Value = 15.02569121839017

$ws.Range("B8").Value = 15.35873794015019
$ws.Range("C8").Value = 7.463938159161071
$ws.Range("D8").Value = 4.476278061671661
$ws.Range("E8").Value = 11.38041901260138
$ws.Range("F8").Value = 60.46433016239353
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 10.82890940980138
$ws.Range("K8").Value = 15.02738111146393

$ws.Range("B9").Value = 15.55252708302371
$ws.Range("C9").Value = 7.498571413083932
$ws.Range("D9").Value = 4.360601037299803
$ws.Range("E9").Value = 11.28510406073295
$ws.Range("F9").Value = 62.85640637024954
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 10.85437582150496
$ws.Range("K9").Value = 15.07457440308378

$ws.Range("B10").Value = 15.72231997927087
$ws.Range("C10").Value = 7.544339504150099
$ws.Range("D10").Value = 4.289035790225291
$ws.Range("E10").Value = 11.23278799735221
$ws.Range("F10").Value = 64.57641039704143
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 10.8807839374306
$ws.Range("K10").Value = 15.135870637658

$ws.Range("B11").Value = 15.80524617340463
$ws.Range("C11").Value = 7.56953299598262
$ws.Range("D11").Value = 4.259491705237365
$ws.Range("E11").Value = 11.21286864452559
$ws.Range("F11").Value = 65.34834109168334
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 10.89445734332756
$ws.Range("K11").Value = 15.16948104229522

$ws.Range("B12").Value = 15.83744016758586
$ws.Range("C12").Value = 7.579696122334142
$ws.Range("D12").Value = 4.248745649394134
$ws.Range("E12").Value = 11.20588613188439
$ws.Range("F12").Value = 65.638947780851
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 10.89987247404844
$ws.Range("K12").Value = 15.18302496355936

$ws.Range("B13").Value = 15.83047184032042
$ws.Range("C13").Value = 7.577479719631865
$ws.Range("D13").Value = 4.251040222347037
$ws.Range("E13").Value = 11.2073649717928
$ws.Range("F13").Value = 65.57643936396347
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 10.89869569909799
$ws.Range("K13").Value = 15.1800718554839

$ws.Range("B14").Value = 15.80787908262963
$ws.Range("C14").Value = 7.570356678564383
$ws.Range("D14").Value = 4.258598719171911
$ws.Range("E14").Value = 11.21228294099428
$ws.Range("F14").Value = 65.37228477515225
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 10.89489810220812
$ws.Range("K14").Value = 15.17057899245953

$ws.Range("B15").Value = 15.79414265725828
$ws.Range("C15").Value = 7.566074525460975
$ws.Range("D15").Value = 4.263286292624492
$ws.Range("E15").Value = 11.21536840989539
$ws.Range("F15").Value = 65.24700626467285
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 10.89260282351217
$ws.Range("K15").Value = 15.16487042503277

$ws.Range("B16").Value = 15.71701279791625
$ws.Range("C16").Value = 7.542780577705485
$ws.Range("D16").Value = 4.29102789241262
$ws.Range("E16").Value = 11.2341680698799
$ws.Range("F16").Value = 64.52573454917501
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 10.87992363812239
$ws.Range("K16").Value = 15.13378877651074

$ws.Range("B17").Value = 15.67113462143466
$ws.Range("C17").Value = 7.529606861266508
$ws.Range("D17").Value = 4.308823789676521
$ws.Range("E17").Value = 11.24669662688739
$ws.Range("F17").Value = 64.08042457533764
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 10.87256963262121
$ws.Range("K17").Value = 15.11618327811755

$ws.Range("B18").Value = 15.64528363626417
$ws.Range("C18").Value = 7.522441789789378
$ws.Range("D18").Value = 4.319342425389291
$ws.Range("E18").Value = 11.2542677307045
$ws.Range("F18").Value = 63.82331621925715
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 10.86849616835249
$ws.Range("K18").Value = 15.10659657147611

$ws.Range("B19").Value = 15.63662391933674
$ws.Range("C19").Value = 7.520086754987458
$ws.Range("D19").Value = 4.322952185414088
$ws.Range("E19").Value = 11.25689378317231
$ws.Range("F19").Value = 63.73610190117288
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 10.86714385974187
$ws.Range("K19").Value = 15.10344353156659

$ws.Range("B20").Value = 15.67596304473768
$ws.Range("C20").Value = 7.530966613262197
$ws.Range("D20").Value = 4.306900044634379
$ws.Range("E20").Value = 11.24532514375554
$ws.Range("F20").Value = 64.12793122869277
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 10.87333630563085
$ws.Range("K20").Value = 15.11800162306205

$ws.Range("B21").Value = 15.81449385483717
$ws.Range("C21").Value = 7.572432038595419
$ws.Range("D21").Value = 4.256366549242182
$ws.Range("E21").Value = 11.21082318189071
$ws.Range("F21").Value = 65.4322977644627
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 10.89600712027858
$ws.Range("K21").Value = 15.17334518286002

$ws.Range("B22").Value = 15.90963275039904
$ws.Range("C22").Value = 7.603158946629282
$ws.Range("D22").Value = 4.225918429153674
$ws.Range("E22").Value = 11.1915423493183
$ws.Range("F22").Value = 66.27474392496447
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 10.91220623311811
$ws.Range("K22").Value = 15.21426972055908

$ws.Range("B23").Value = 15.85844340276871
$ws.Range("C23").Value = 7.586429939292731
$ws.Range("D23").Value = 4.241930321830986
$ws.Range("E23").Value = 11.20153302806633
$ws.Range("F23").Value = 65.82609501794721
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 10.90343448904541
$ws.Range("K23").Value = 15.19199518623194

$ws.Range("B24").Value = 15.67377847738108
$ws.Range("C24").Value = 7.53035059602268
$ws.Range("D24").Value = 4.3077688748034
$ws.Range("E24").Value = 11.24594404440493
$ws.Range("F24").Value = 64.106456862469
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 10.87298921157522
$ws.Range("K24").Value = 15.11717788260918

$ws.Range("B25").Value = 15.49520254664401
$ws.Range("C25").Value = 7.485624713897897
$ws.Range("D25").Value = 4.389569891708596
$ws.Range("E25").Value = 11.30779105151969
$ws.Range("F25").Value = 62.21502169856142
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 10.84613109414047
$ws.Range("K25").Value = 15.05711738057904
